$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.10125
$ws.Range("E2").Value = -0.112
$ws.Range("I2").Value = 0.001809478037358186
$ws.Range("J2").Value = 0.001809478037358186
$ws.Range("K2").Value = -0.09000000000000008
$ws.Range("L2").Value = -0.006122448979591842
$ws.Range("U2").Value = 47.2
$ws.Range("V2").Value = 0.214838416021848
$ws.Range("W2").Value = 0.005828877005347594
$ws.Range("X2").Value = 0.03215313721310308
$ws.Range("Y2").Value = -0.02632426020775548
$ws.Range("Z2").Value = 0.416053396362805
$ws.Range("AA2").Value = 0.0004231249093625263
$ws.Range("AB2").Value = 0.03198870164868391
$ws.Range("AC2").Value = -0.03156557673932138
$ws.Range("AD2").Value = 5.78
$ws.Range("AE2").Value = 0.0320033642541733
$ws.Range("AF2").Value = 5.812003364254173
$ws.Range("AG2").Value = -41.38799663574583
$ws.Range("AH2").Value = 0.02577247897029428
$ws.Range("AI2").Value = 0.05719799995891788
$ws.Range("AJ2").Value = -0.2321099861751808
$ws.Range("AK2").Value = -0.7606409262066534
$ws.Range("AN2").Value = 175.1515151515152
$ws.Range("AP2").Value = -1254.181716234722

# Row 3
$ws.Range("D3").Value = -0.122
$ws.Range("I3").Value = 0.002582458946520907
$ws.Range("J3").Value = 0.002582458946520907
$ws.Range("K3").Value = -1.59
$ws.Range("L3").Value = -0.1543689320388349
$ws.Range("U3").Value = 27.6
$ws.Range("V3").Value = 0.1970021413276231
$ws.Range("W3").Value = -0.02834224598930481
$ws.Range("X3").Value = 0.03143236596107401
$ws.Range("Y3").Value = -0.05977461195037882
$ws.Range("Z3").Value = 0.3276914894872276
$ws.Range("AA3").Value = 0.0008462498187250527
$ws.Range("AB3").Value = 0.03143158834223347
$ws.Range("AC3").Value = -0.03058533852350842
$ws.Range("AE3").Value = 0.0320033642541733
$ws.Range("AF3").Value = 0.0320033642541733
$ws.Range("AG3").Value = -27.56799663574583
$ws.Range("AH3").Value = 0.0002283801236394579
$ws.Range("AI3").Value = 0.0005621331125380214
$ws.Range("AJ3").Value = -0.2449791686948925
$ws.Range("AK3").Value = -0.9398606802746358
$ws.Range("AP3").Value = -835.3938374468432

# Row 4
$ws.Range("D4").Value = -0.0805
$ws.Range("E4").Value = -0.112
$ws.Range("K4").Value = 1.5
$ws.Range("L4").Value = 0.3409090909090909
$ws.Range("U4").Value = 19.6
$ws.Range("V4").Value = 0.2462311557788945
$ws.Range("W4").Value = 0.04
$ws.Range("X4").Value = 0.03287390846513215
$ws.Range("Y4").Value = 0.007126091534867852
$ws.Range("Z4").Value = 1.128205128205129
$ws.Range("AB4").Value = 0.03254581495513435
$ws.Range("AC4").Value = -0.03254581495513435
$ws.Range("AD4").Value = 5.78
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 5.78
$ws.Range("AG4").Value = -13.82
$ws.Range("AH4").Value = 0.06769735301007263
$ws.Range("AI4").Value = 0.1293643688451209
$ws.Range("AJ4").Value = -0.2100942535725144
$ws.Range("AK4").Value = -0.55103668261563
